$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- Row 65: new command "添加列车径路" (AddTrainPath) ---
$ws.Cells.Item(65, 1).Value = 63
$ws.Cells.Item(65, 2).Value = "添加列车径路"
$ws.Cells.Item(65, 3).Value = "AddTrainPath"
$ws.Cells.Item(65, 4).Value = "支持"
$ws.Cells.Item(65, 5).Value = "否"
$ws.Cells.Item(65, 6).Value = 45150

# --- Row 66: new command "删除列车径路" (RemoveTrainPath) ---
$ws.Cells.Item(66, 1).Value = 64
$ws.Cells.Item(66, 2).Value = "删除列车径路"
$ws.Cells.Item(66, 3).Value = "RemoveTrainPath"
$ws.Cells.Item(66, 4).Value = "支持"
$ws.Cells.Item(66, 5).Value = "否"
$ws.Cells.Item(66, 6).Value = 45150

# Match the date display format used by the rest of column F ("完成时间")
$ws.Range("F65:F66").NumberFormat = "yyyy""年""m""月""d""日"";@"

# Update the sheet's active selection to reflect where editing ended up
$ws.Range("C69").Select() | Out-Null
